$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.01725120502155203
$ws.Range("C2").Value = 1.384110966838059
$ws.Range("D2").Value = 8.376137918165924
$ws.Range("E2").Value = 2.894155821334768
$ws.Range("F2").Value = 2.959148581664683
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = -0.1059444359486
$ws.Range("C3").Value = 1.363387133547353
$ws.Range("D3").Value = 7.96756717039053
$ws.Range("E3").Value = 2.822687933582197
$ws.Range("F3").Value = 2.887077487796805
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = -0.5314575476434735
$ws.Range("C4").Value = 0.935032954889049
$ws.Range("D4").Value = 3.962372580821567
$ws.Range("E4").Value = 1.990570918310013
$ws.Range("F4").Value = 1.965686070731558
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = -0.1228908575652036
$ws.Range("C5").Value = 0.5948419731728133
$ws.Range("D5").Value = 0.8938785104026519
$ws.Range("E5").Value = 0.9454514849544908
$ws.Range("F5").Value = 0.9617836213211481
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = -0.07803016097268373
$ws.Range("C6").Value = 0.595167622792783
$ws.Range("D6").Value = 0.6624405579447651
$ws.Range("E6").Value = 0.8139045140216173
$ws.Range("F6").Value = 0.8323555993064871
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = -0.1177682127167521
$ws.Range("C7").Value = 0.5278966303169981
$ws.Range("D7").Value = 0.6234154974736094
$ws.Range("E7").Value = 0.7895666516980118
$ws.Range("F7").Value = 0.8033690317453926
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = 0.008373115946536357
$ws.Range("C8").Value = 0.403612244087397
$ws.Range("D8").Value = 0.3646454851063948
$ws.Range("E8").Value = 0.6038588287889768
$ws.Range("F8").Value = 0.6223835931625885
$ws.Range("G8").Value = 17

$ws.Range("B9").Value = 0.005946962572950867
$ws.Range("C9").Value = 0.3740864020203869
$ws.Range("D9").Value = 0.2564549243757536
$ws.Range("E9").Value = 0.5064137877030538
$ws.Range("F9").Value = 0.5229858460921323
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = 0.01722052028762557
$ws.Range("C10").Value = 0.3408956206005759
$ws.Range("D10").Value = 0.2330792879415247
$ws.Range("E10").Value = 0.4827828579615527
$ws.Range("F10").Value = 0.4994097318934355
$ws.Range("G10").Value = 15

$ws.Range("B11").Value = 0.04981916865326922
$ws.Range("C11").Value = 0.3023652478618333
$ws.Range("D11").Value = 0.214504976889282
$ws.Range("E11").Value = 0.4631468200142175
$ws.Range("F11").Value = 0.4778414914637334
$ws.Range("G11").Value = 14
